$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-17 02:51:05"
$wsZh.Range("H4").Value = "2016-03-17 02:52:00"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-17 02:51:13"
$wsDe.Range("H4").Value = "2016-03-17 02:52:13"
